$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the two "Incomplete ||" status cells to reflect the newly completed work
$ws.Range("F7").Value = "Completed || 03/02 - 04/02"
$ws.Range("F12").Value = "Completed || 03/02 - 04/02"

# Update the last active selection to match the saved state
$ws.Range("E19").Select()
